$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows sourced from DGS's 2021/12/22 and 2021/12/24 reports.
# (row 123 -> 2021/12/22, row 124 -> 2021/12/24)
$newRows = @(
    @{ Row = 123; Date = "2021/12/22"; B = 579.3; C = 582.3; D = 1.07; E = 1.07 },
    @{ Row = 124; Date = "2021/12/24"; B = 630.8; C = 633.1; D = 1.11; E = 1.11 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Date
    $cellA.NumberFormat = "yyyy/mm/dd"

    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellB.NumberFormat = "0.00"
    $cellB.Value = $r.B

    $cellC = $ws.Cells.Item($rowNum, 3)
    $cellC.NumberFormat = "0.00"
    $cellC.Value = $r.C

    $cellD = $ws.Cells.Item($rowNum, 4)
    $cellD.NumberFormat = "0.00"
    $cellD.Value = $r.D

    $cellE = $ws.Cells.Item($rowNum, 5)
    $cellE.NumberFormat = "0.00"
    $cellE.Value = $r.E
}

# Update the active cell to reflect the appended rows, matching the
# author's cursor position after entering the new data.
$ws.Range("A125").Select()
